$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 2121
$ws.Range("I111").Value = 2128.9285
$ws.Range("J111").Value = 2109.9
$ws.Range("K111").Value = 6386.7855
$ws.Range("L111").Value = 6329.700000000001
$ws.Range("M111").Value = -3319.7855
$ws.Range("N111").Value = -12463.7

$ws.Range("H116").Value = 8698288
$ws.Range("I116").Value = 18183830
$ws.Range("J116").Value = 3208.25
$ws.Range("K116").Value = 18183830
$ws.Range("L116").Value = 3208.25
$ws.Range("M116").Value = -18180388
$ws.Range("N116").Value = -10092.25

$ws.Range("H137").Value = 1329.8518
$ws.Range("I137").Value = 1155.6666
$ws.Range("J137").Value = 1678.2222
$ws.Range("K137").Value = 3466.9998
$ws.Range("L137").Value = 5034.6666
$ws.Range("M137").Value = -916.9998000000001
$ws.Range("N137").Value = -10134.6666

$ws.Range("H138").Value = 4195.75
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 4195.75
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 12587.25
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -22867.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 46800
$ws.Range("I2").Value = 1278.1
$ws.Range("J2").Value = 84734.914
$ws.Range("K2").Value = 1278.1
$ws.Range("L2").Value = 84734.914
$ws.Range("M2").Value = -1165.1
$ws.Range("N2").Value = -84960.914

$ws.Range("H61").Value = 1365.0613
$ws.Range("I61").Value = 1150.2439
$ws.Range("J61").Value = 2466
$ws.Range("K61").Value = 1150.2439
$ws.Range("L61").Value = 2466
$ws.Range("M61").Value = -938.2438999999999
$ws.Range("N61").Value = -2890

$ws.Range("H97").Value = 1786.5714
$ws.Range("I97").Value = 1474.6364
$ws.Range("J97").Value = 2930.3333
$ws.Range("K97").Value = 1474.6364
$ws.Range("L97").Value = 2930.3333
$ws.Range("M97").Value = -978.6364000000001
$ws.Range("N97").Value = -3922.3333

$ws.Range("H101").Value = 53361.2
$ws.Range("J101").Value = 53361.2
$ws.Range("L101").Value = 53361.2
$ws.Range("N101").Value = -59851.2

$ws.Range("H106").Value = 50000
$ws.Range("J106").Value = 50000
$ws.Range("L106").Value = 50000
$ws.Range("N106").Value = -52524

$ws.Range("H110").Value = 1334.6154
$ws.Range("I110").Value = 1343.9
$ws.Range("K110").Value = 1343.9
$ws.Range("M110").Value = 701.0999999999999

$ws.Range("H116").Value = 46800
$ws.Range("I116").Value = 1278.1
$ws.Range("J116").Value = 84734.914
$ws.Range("K116").Value = 1278.1
$ws.Range("L116").Value = 84734.914
$ws.Range("M116").Value = 1015.9
$ws.Range("N116").Value = -89322.914

$ws.Range("H122").Value = 1674.8572
$ws.Range("I122").Value = 1554
$ws.Range("J122").Value = 2400
$ws.Range("K122").Value = 4662
$ws.Range("L122").Value = 7200
$ws.Range("M122").Value = -2212
$ws.Range("N122").Value = -12100

$ws.Range("H123").Value = 24313.166
$ws.Range("J123").Value = 24313.166
$ws.Range("L123").Value = 24313.166
$ws.Range("N123").Value = -34113.166

$ws.Range("H132").Value = 1913.7122
$ws.Range("I132").Value = 1497.5957
$ws.Range("J132").Value = 2943.0527
$ws.Range("K132").Value = 4492.7871
$ws.Range("L132").Value = 8829.158100000001
$ws.Range("M132").Value = -1962.7871
$ws.Range("N132").Value = -13889.1581

$ws.Range("H136").Value = 1365.0613
$ws.Range("I136").Value = 1150.2439
$ws.Range("J136").Value = 2466
$ws.Range("K136").Value = 3450.7317
$ws.Range("L136").Value = 7398
$ws.Range("M136").Value = -900.7316999999998
$ws.Range("N136").Value = -12498

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 46800
$ws.Range("I3").Value = 1278.1
$ws.Range("J3").Value = 84734.914
$ws.Range("K3").Value = 1278.1
$ws.Range("L3").Value = 84734.914
$ws.Range("M3").Value = -1164.1
$ws.Range("N3").Value = -84962.914

$ws.Range("H109").Value = 20416.188
$ws.Range("J109").Value = 20416.188
$ws.Range("L109").Value = 20416.188
$ws.Range("N109").Value = -23190.188

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H100").Value = 50000
$ws.Range("J100").Value = 50000
$ws.Range("L100").Value = 50000
$ws.Range("N100").Value = -52164

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1110.1818
$ws.Range("I2").Value = 2124.4
$ws.Range("J2").Value = 265
$ws.Range("K2").Value = 12746.4
$ws.Range("L2").Value = 1590
$ws.Range("M2").Value = -12633.4
$ws.Range("N2").Value = -1816

$ws.Range("H38").Value = 129.25
$ws.Range("I38").Value = 109.666664
$ws.Range("J38").Value = 188
$ws.Range("K38").Value = 328.999992
$ws.Range("L38").Value = 564
$ws.Range("M38").Value = 18.00000799999998
$ws.Range("N38").Value = -1258

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 21999.75
$ws.Range("J42").Value = 21999.75
$ws.Range("L42").Value = 21999.75
$ws.Range("N42").Value = -22969.75

$ws.Range("H45").Value = 38659
$ws.Range("J45").Value = 38659
$ws.Range("L45").Value = 38659
$ws.Range("N45").Value = -39777

$ws.Range("H51").Value = 34582.668
$ws.Range("J51").Value = 34582.668
$ws.Range("L51").Value = 34582.668
$ws.Range("N51").Value = -35600.668

$ws.Range("H97").Value = 126725.086
$ws.Range("I97").Value = 64548.75
$ws.Range("J97").Value = 251077.75
$ws.Range("K97").Value = 64548.75
$ws.Range("L97").Value = 251077.75
$ws.Range("M97").Value = -64052.75
$ws.Range("N97").Value = -252069.75

$ws.Range("H102").Value = 5166.6665
$ws.Range("I102").Value = 4333.3335
$ws.Range("J102").Value = 5583.3335
$ws.Range("K102").Value = 4333.3335
$ws.Range("L102").Value = 5583.3335
$ws.Range("M102").Value = -2711.3335
$ws.Range("N102").Value = -8827.333500000001

$ws.Range("H109").Value = 13882.6
$ws.Range("J109").Value = 13882.6
$ws.Range("L109").Value = 13882.6
$ws.Range("N109").Value = -15962.6

$ws.Range("H113").Value = 1885.5
$ws.Range("I113").Value = 1340
$ws.Range("J113").Value = 2275.1428
$ws.Range("K113").Value = 1340
$ws.Range("L113").Value = 2275.1428
$ws.Range("M113").Value = 830
$ws.Range("N113").Value = -6615.1428

$ws.Range("H115").Value = 21999.75
$ws.Range("J115").Value = 21999.75
$ws.Range("L115").Value = 21999.75
$ws.Range("N115").Value = -24349.75

$ws.Range("H122").Value = 3966.75
$ws.Range("I122").Value = 3007
$ws.Range("J122").Value = 4286.6665
$ws.Range("K122").Value = 9021
$ws.Range("L122").Value = 12859.9995
$ws.Range("M122").Value = -6571
$ws.Range("N122").Value = -17759.9995

$ws.Range("H123").Value = 10553.389
$ws.Range("J123").Value = 10553.389
$ws.Range("L123").Value = 10553.389
$ws.Range("N123").Value = -15453.389

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 6126.154
$ws.Range("I100").Value = 8677.143
$ws.Range("J100").Value = 3150
$ws.Range("K100").Value = 8677.143
$ws.Range("L100").Value = 3150
$ws.Range("M100").Value = -8136.143
$ws.Range("N100").Value = -4232

$ws.Range("H116").Value = 48680
$ws.Range("J116").Value = 48680
$ws.Range("L116").Value = 48680
$ws.Range("N116").Value = -57858

$ws.Range("H122").Value = 22506216
$ws.Range("I122").Value = 27783756
$ws.Range("J122").Value = 18188228
$ws.Range("K122").Value = 83351268
$ws.Range("L122").Value = 54564684
$ws.Range("M122").Value = -83348818
$ws.Range("N122").Value = -54569584

$ws.Range("H131").Value = 59326
$ws.Range("J131").Value = 59326
$ws.Range("L131").Value = 59326
$ws.Range("N131").Value = -69406

$ws.Range("H134").Value = 47143
$ws.Range("J134").Value = 47143
$ws.Range("L134").Value = 47143
$ws.Range("N134").Value = -57283

$ws.Range("H136").Value = 1893.6875
$ws.Range("I136").Value = 1590.7322
$ws.Range("J136").Value = 4014.375
$ws.Range("K136").Value = 4772.196599999999
$ws.Range("L136").Value = 12043.125
$ws.Range("M136").Value = -2222.196599999999
$ws.Range("N136").Value = -17143.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 22151.666
$ws.Range("J123").Value = 22151.666
$ws.Range("L123").Value = 22151.666
$ws.Range("N123").Value = -31951.666

$ws.Range("H126").Value = 8186.353
$ws.Range("I126").Value = 9018.200000000001
$ws.Range("J126").Value = 1947.5
$ws.Range("K126").Value = 27054.6
$ws.Range("L126").Value = 5842.5
$ws.Range("M126").Value = -24584.6
$ws.Range("N126").Value = -10782.5

$ws.Range("H136").Value = 1223.3871
$ws.Range("I136").Value = 1311.4286
$ws.Range("J136").Value = 401.66666
$ws.Range("K136").Value = 3934.2858
$ws.Range("L136").Value = 1204.99998
$ws.Range("M136").Value = -1384.2858
$ws.Range("N136").Value = -6304.999980000001
